$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = " 67"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = 220
$ws.Range("D2").Value = 216
$ws.Range("E2").Value = 1.25
$ws.Range("F2").Value = 149.7
$ws.Range("G2").Value = 0.74
$ws.Range("H2").Value = 0.78
$ws.Range("J2").Value = 0.1
$ws.Range("L2").Value = 171
$ws.Range("M2").Value = 137
$ws.Range("N2").Value = 35
$ws.Range("O2").Value = 23
$ws.Range("P2").Value = 17

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = " 27"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = 87
$ws.Range("C3").Value = 0.76
$ws.Range("D3").Value = 179.8
$ws.Range("E3").Value = 0.86
$ws.Range("F3").Value = 121.5
$ws.Range("H3").Value = 0.68
$ws.Range("I3").Value = 0.11
$ws.Range("J3").Value = 0.08
$ws.Range("K3").Value = 0.06
$ws.Range("L3").Value = 59
$ws.Range("M3").Value = 69
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 7
$ws.Range("P3").Value = 5

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = " 7"
$ws.Range("A4").Style = "Normal"
